$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")

$ws.Cells.Item(2, 1).Value = "Última actualización: 10:30:21"
$ws.Cells.Item(3, 1).Value = "Total filas: 148"
$ws.Cells.Item(32, 1).Value = "06:44:15"
$ws.Cells.Item(32, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(32, 4).Value = 21
$ws.Cells.Item(33, 1).Value = "05:31:23"
$ws.Cells.Item(33, 3).Value = "15_ABASTO"
$ws.Cells.Item(33, 4).Value = 94
$ws.Cells.Item(96, 1).Value = "07:57:27"
$ws.Cells.Item(96, 3).Value = "17_ROMERO"
$ws.Cells.Item(96, 4).Value = 86
$ws.Cells.Item(97, 1).Value = "07:31:43"
$ws.Cells.Item(97, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(97, 4).Value = 112
$ws.Cells.Item(104, 1).Value = "07:57:27"
$ws.Cells.Item(104, 3).Value = "215C_EL PATO"
$ws.Cells.Item(104, 4).Value = 105
$ws.Cells.Item(105, 1).Value = "09:31:25"
$ws.Cells.Item(105, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(105, 4).Value = 11
$ws.Cells.Item(120, 1).Value = "10:30:21"
$ws.Cells.Item(120, 2).Value = "10:31"
$ws.Cells.Item(120, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(120, 4).Value = 1
$ws.Cells.Item(121, 1).Value = "10:30:21"
$ws.Cells.Item(121, 2).Value = "10:34"
$ws.Cells.Item(121, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(121, 4).Value = 4
$ws.Cells.Item(122, 2).Value = "10:41"
$ws.Cells.Item(122, 3).Value = "17_ROMERO"
$ws.Cells.Item(122, 4).Value = 106
$ws.Cells.Item(123, 2).Value = "10:42"
$ws.Cells.Item(123, 3).Value = "17_ROMERO"
$ws.Cells.Item(123, 4).Value = 71
$ws.Cells.Item(124, 1).Value = "08:55:44"
$ws.Cells.Item(124, 2).Value = "10:43"
$ws.Cells.Item(124, 3).Value = "14_ABASTO"
$ws.Cells.Item(124, 4).Value = 108
$ws.Cells.Item(125, 1).Value = "10:30:21"
$ws.Cells.Item(125, 2).Value = "10:46"
$ws.Cells.Item(125, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(125, 4).Value = 16
$ws.Cells.Item(126, 1).Value = "10:30:21"
$ws.Cells.Item(126, 2).Value = "10:52"
$ws.Cells.Item(126, 3).Value = "15_ABASTO"
$ws.Cells.Item(126, 4).Value = 22
$ws.Cells.Item(127, 1).Value = "10:30:21"
$ws.Cells.Item(127, 2).Value = "10:53"
$ws.Cells.Item(127, 3).Value = "10_OLMOS"
$ws.Cells.Item(127, 4).Value = 23
$ws.Cells.Item(128, 1).Value = "10:30:21"
$ws.Cells.Item(128, 2).Value = "10:57"
$ws.Cells.Item(128, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(128, 4).Value = 27
$ws.Cells.Item(129, 2).Value = "10:59"
$ws.Cells.Item(129, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(129, 4).Value = 88
$ws.Cells.Item(130, 1).Value = "10:30:21"
$ws.Cells.Item(130, 2).Value = "10:59"
$ws.Cells.Item(130, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(130, 4).Value = 29
$ws.Cells.Item(130, 5).Value = "LP1912"
$ws.Cells.Item(131, 1).Value = "09:31:25"
$ws.Cells.Item(131, 2).Value = "11:02"
$ws.Cells.Item(131, 3).Value = "215C_EL PATO"
$ws.Cells.Item(131, 4).Value = 91
$ws.Cells.Item(131, 5).Value = "LP1912"
$ws.Cells.Item(132, 1).Value = "10:30:21"
$ws.Cells.Item(132, 2).Value = "11:03"
$ws.Cells.Item(132, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(132, 4).Value = 33
$ws.Cells.Item(132, 5).Value = "LP1912"
$ws.Cells.Item(133, 1).Value = "10:30:21"
$ws.Cells.Item(133, 2).Value = "11:06"
$ws.Cells.Item(133, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(133, 4).Value = 36
$ws.Cells.Item(133, 5).Value = "LP1912"
$ws.Cells.Item(134, 1).Value = "10:30:21"
$ws.Cells.Item(134, 2).Value = "11:11"
$ws.Cells.Item(134, 3).Value = "10_OLMOS"
$ws.Cells.Item(134, 4).Value = 41
$ws.Cells.Item(134, 5).Value = "LP1912"
$ws.Cells.Item(135, 1).Value = "09:31:25"
$ws.Cells.Item(135, 2).Value = "11:17"
$ws.Cells.Item(135, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(135, 4).Value = 106
$ws.Cells.Item(135, 5).Value = "LP1912"
$ws.Cells.Item(136, 1).Value = "09:31:25"
$ws.Cells.Item(136, 2).Value = "11:19"
$ws.Cells.Item(136, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(136, 4).Value = 108
$ws.Cells.Item(136, 5).Value = "LP1912"
$ws.Cells.Item(137, 1).Value = "09:31:25"
$ws.Cells.Item(137, 2).Value = "11:21"
$ws.Cells.Item(137, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(137, 4).Value = 110
$ws.Cells.Item(137, 5).Value = "LP1912"
$ws.Cells.Item(138, 1).Value = "09:31:25"
$ws.Cells.Item(138, 2).Value = "11:26"
$ws.Cells.Item(138, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(138, 4).Value = 115
$ws.Cells.Item(138, 5).Value = "LP1912"
$ws.Cells.Item(139, 1).Value = "09:31:25"
$ws.Cells.Item(139, 2).Value = "11:27"
$ws.Cells.Item(139, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(139, 4).Value = 116
$ws.Cells.Item(139, 5).Value = "LP1912"
$ws.Cells.Item(140, 1).Value = "10:30:21"
$ws.Cells.Item(140, 2).Value = "11:32"
$ws.Cells.Item(140, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(140, 4).Value = 62
$ws.Cells.Item(140, 5).Value = "LP1912"
$ws.Cells.Item(141, 1).Value = "10:30:21"
$ws.Cells.Item(141, 2).Value = "11:35"
$ws.Cells.Item(141, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(141, 4).Value = 65
$ws.Cells.Item(141, 5).Value = "LP1912"
$ws.Cells.Item(142, 1).Value = "10:30:21"
$ws.Cells.Item(142, 2).Value = "11:39"
$ws.Cells.Item(142, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(142, 4).Value = 69
$ws.Cells.Item(142, 5).Value = "LP1912"
$ws.Cells.Item(143, 1).Value = "10:30:21"
$ws.Cells.Item(143, 2).Value = "11:42"
$ws.Cells.Item(143, 3).Value = "17_ROMERO"
$ws.Cells.Item(143, 4).Value = 72
$ws.Cells.Item(143, 5).Value = "LP1912"
$ws.Cells.Item(144, 1).Value = "10:30:21"
$ws.Cells.Item(144, 2).Value = "11:48"
$ws.Cells.Item(144, 3).Value = "10_OLMOS"
$ws.Cells.Item(144, 4).Value = 78
$ws.Cells.Item(144, 5).Value = "LP1912"
$ws.Cells.Item(145, 1).Value = "10:30:21"
$ws.Cells.Item(145, 2).Value = "11:51"
$ws.Cells.Item(145, 3).Value = "215B_EL PATO"
$ws.Cells.Item(145, 4).Value = 81
$ws.Cells.Item(145, 5).Value = "LP1912"
$ws.Cells.Item(146, 1).Value = "10:30:21"
$ws.Cells.Item(146, 2).Value = "11:54"
$ws.Cells.Item(146, 3).Value = "15_ABASTO"
$ws.Cells.Item(146, 4).Value = 84
$ws.Cells.Item(146, 5).Value = "LP1912"
$ws.Cells.Item(147, 1).Value = "10:30:21"
$ws.Cells.Item(147, 2).Value = "11:59"
$ws.Cells.Item(147, 3).Value = "225_GOMEZ"
$ws.Cells.Item(147, 4).Value = 89
$ws.Cells.Item(147, 5).Value = "LP1912"
$ws.Cells.Item(148, 1).Value = "10:30:21"
$ws.Cells.Item(148, 2).Value = "12:02"
$ws.Cells.Item(148, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(148, 4).Value = 92
$ws.Cells.Item(148, 5).Value = "LP1912"
$ws.Cells.Item(149, 1).Value = "10:30:21"
$ws.Cells.Item(149, 2).Value = "12:06"
$ws.Cells.Item(149, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(149, 4).Value = 96
$ws.Cells.Item(149, 5).Value = "LP1912"
$ws.Cells.Item(150, 1).Value = "10:30:21"
$ws.Cells.Item(150, 2).Value = "12:14"
$ws.Cells.Item(150, 3).Value = "17_ROMERO"
$ws.Cells.Item(150, 4).Value = 104
$ws.Cells.Item(150, 5).Value = "LP1912"
$ws.Cells.Item(151, 1).Value = "10:30:21"
$ws.Cells.Item(151, 2).Value = "12:17"
$ws.Cells.Item(151, 3).Value = "14_ABASTO"
$ws.Cells.Item(151, 4).Value = 107
$ws.Cells.Item(151, 5).Value = "LP1912"
$ws.Cells.Item(152, 1).Value = "10:30:21"
$ws.Cells.Item(152, 2).Value = "12:20"
$ws.Cells.Item(152, 3).Value = "215A_EL PATO"
$ws.Cells.Item(152, 4).Value = 110
$ws.Cells.Item(152, 5).Value = "LP1912"
$ws.Cells.Item(153, 1).Value = "10:30:21"
$ws.Cells.Item(153, 2).Value = "12:21"
$ws.Cells.Item(153, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(153, 4).Value = 111
$ws.Cells.Item(153, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2, 1).Value = "Última actualización: 10:30:21"
$ws.Cells.Item(3, 1).Value = "Total filas: 17"
$ws.Cells.Item(21, 1).Value = "10:30:21"
$ws.Cells.Item(21, 2).Value = "11:51"
$ws.Cells.Item(21, 3).Value = "215B_EL PATO"
$ws.Cells.Item(21, 4).Value = 81
$ws.Cells.Item(21, 5).Value = "LP1912"
$ws.Cells.Item(22, 1).Value = "10:30:21"
$ws.Cells.Item(22, 2).Value = "12:20"
$ws.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws.Cells.Item(22, 4).Value = 110
$ws.Cells.Item(22, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2, 1).Value = "Última actualización: 10:30:21"
$ws.Cells.Item(3, 1).Value = "Total filas: 25"
$ws.Cells.Item(30, 1).Value = "10:30:21"
$ws.Cells.Item(30, 2).Value = "12:04"
$ws.Cells.Item(30, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(30, 4).Value = 94
$ws.Cells.Item(30, 5).Value = "L6173"
